$wb = $excel.ActiveWorkbook

$values = @{
    "C11" = 0.676738050782717
    "D11" = -0.2999999999999998
    "E11" = 0.6333398306071416
    "F11" = -0.038999999999999924
    "G11" = 1.5829618029997903
    "H11" = 16.12947350163202
    "I11" = 2.03501970083987
}

foreach ($sheetName in @("Test 1", "Test 2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $values.Keys) {
        $ws.Range($addr).Value = $values[$addr]
    }
}
